$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the date-like text columns (A, D, E) as Text so Excel
# does not auto-convert strings like "2024-04-24" into date serials.
# (Each column range is set separately -- a single multi-area Range
# does not reliably propagate NumberFormat to every area.)
$colA = $ws.Range("A2:A14")
$colD = $ws.Range("D2:D14")
$colE = $ws.Range("E2:E14")
$colA.NumberFormat = "@"
$colD.NumberFormat = "@"
$colE.NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '2024-04-24'
$ws.Range("B2").Value = '코칩'
$ws.Range("C2").Value = '한국'
$ws.Range("D2").Value = '2024-04-29'
$ws.Range("E2").Value = '2024-05-07'
$ws.Range("F2").Value = 27000000
$ws.Range("G2").Value = 1500000
$ws.Range("H2").Value = '-'
$ws.Range("I2").Value = 11000
$ws.Range("J2").Value = 14000
$ws.Range("K2").Value = '-'
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = '-'
$ws.Range("N2").Value = '-'
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = '-'
$ws.Range("Q2").Value = '-'
$ws.Range("R2").Value = '734.49 : 1'
$ws.Range("S2").Value = '-'
$ws.Range("T2").Value = '-'

# Row 3
$ws.Range("A3").Value = '2024-04-23'
$ws.Range("B3").Value = 'SK증권제12호스팩'
$ws.Range("C3").Value = 'SK'
$ws.Range("D3").Value = '2024-04-26'
$ws.Range("E3").Value = '2024-05-07'
$ws.Range("F3").Value = 6000000
$ws.Range("G3").Value = 3000000
$ws.Range("H3").Value = '-'
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = '-'
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = '-'
$ws.Range("N3").Value = '-'
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = '-'
$ws.Range("Q3").Value = '-'
$ws.Range("R3").Value = '1091.48 : 1'
$ws.Range("S3").Value = '-'
$ws.Range("T3").Value = '-'

# Row 4
$ws.Range("A4").Value = '2024-04-23'
$ws.Range("B4").Value = '민테크'
$ws.Range("C4").Value = 'KB'
$ws.Range("D4").Value = '2024-04-26'
$ws.Range("E4").Value = '2024-05-03'
$ws.Range("F4").Value = 31500000
$ws.Range("G4").Value = 3000000
$ws.Range("H4").Value = '-'
$ws.Range("I4").Value = 6500
$ws.Range("J4").Value = 8500
$ws.Range("K4").Value = '-'
$ws.Range("L4").Value = 10500
$ws.Range("M4").Value = '-'
$ws.Range("N4").Value = '-'
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = '-'
$ws.Range("Q4").Value = '-'
$ws.Range("R4").Value = '1529.43 : 1'
$ws.Range("S4").Value = '-'
$ws.Range("T4").Value = '-'

# Row 5
$ws.Range("A5").Value = '2024-04-22'
$ws.Range("B5").Value = '디앤디파마텍'
$ws.Range("C5").Value = '한국'
$ws.Range("D5").Value = '2024-04-25'
$ws.Range("E5").Value = '2024-05-02'
$ws.Range("F5").Value = 36300000
$ws.Range("G5").Value = 1100000
$ws.Range("H5").Value = '-'
$ws.Range("I5").Value = 22000
$ws.Range("J5").Value = 26000
$ws.Range("K5").Value = '-'
$ws.Range("L5").Value = 33000
$ws.Range("M5").Value = '-'
$ws.Range("N5").Value = '-'
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = '-'
$ws.Range("Q5").Value = '-'
$ws.Range("R5").Value = '1544 : 1'
$ws.Range("S5").Value = '-'
$ws.Range("T5").Value = '-'

# Row 6
$ws.Range("A6").Value = '2024-04-22'
$ws.Range("B6").Value = '유안타제16호스팩'
$ws.Range("C6").Value = '유안타'
$ws.Range("D6").Value = '2024-04-25'
$ws.Range("E6").Value = '2024-05-02'
$ws.Range("F6").Value = 10300000
$ws.Range("G6").Value = 5150000
$ws.Range("H6").Value = '-'
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = '-'
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = '-'
$ws.Range("N6").Value = '-'
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = '-'
$ws.Range("Q6").Value = '-'
$ws.Range("R6").Value = '334 : 1'
$ws.Range("S6").Value = '-'
$ws.Range("T6").Value = '-'

# Row 7
$ws.Range("A7").Value = '2024-04-18'
$ws.Range("B7").Value = '제일엠앤에스'
$ws.Range("C7").Value = 'KB'
$ws.Range("D7").Value = '2024-04-23'
$ws.Range("E7").Value = '2024-04-30'
$ws.Range("F7").Value = 52800000
$ws.Range("G7").Value = 2400000
$ws.Range("H7").Value = '-'
$ws.Range("I7").Value = 15000
$ws.Range("J7").Value = 18000
$ws.Range("K7").Value = '-'
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = '-'
$ws.Range("N7").Value = '-'
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = '-'
$ws.Range("Q7").Value = '-'
$ws.Range("R7").Value = '1438.96 : 1'
$ws.Range("S7").Value = '-'
$ws.Range("T7").Value = '-'

# Row 8
$ws.Range("A8").Value = '2024-04-15'
$ws.Range("B8").Value = '하나33호스팩'
$ws.Range("C8").Value = '하나'
$ws.Range("D8").Value = '2024-04-18'
$ws.Range("E8").Value = '2024-04-24'
$ws.Range("F8").Value = 7000000
$ws.Range("G8").Value = 3500000
$ws.Range("H8").Value = '-'
$ws.Range("I8").Value = 2000
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = '-'
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = '-'
$ws.Range("N8").Value = '-'
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = '-'
$ws.Range("Q8").Value = '-'
$ws.Range("R8").Value = '2248.41 : 1'
$ws.Range("S8").Value = '-'
$ws.Range("T8").Value = '-'

# Row 9
$ws.Range("A9").Value = '2024-04-11'
$ws.Range("B9").Value = '신한제13호스팩'
$ws.Range("C9").Value = '신한'
$ws.Range("D9").Value = '2024-04-15'
$ws.Range("E9").Value = '2024-04-22'
$ws.Range("F9").Value = 6000000
$ws.Range("G9").Value = 3000000
$ws.Range("H9").Value = '-'
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = '-'
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = '-'
$ws.Range("N9").Value = '-'
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = '-'
$ws.Range("Q9").Value = '-'
$ws.Range("R9").Value = '1337.88 : 1'
$ws.Range("S9").Value = '-'
$ws.Range("T9").Value = '-'

# Row 10
$ws.Range("A10").Value = '2024-04-02'
$ws.Range("B10").Value = '신한제12호스팩'
$ws.Range("C10").Value = '신한'
$ws.Range("D10").Value = '2024-04-05'
$ws.Range("E10").Value = '2024-04-15'
$ws.Range("F10").Value = 10000000
$ws.Range("G10").Value = 5000000
$ws.Range("H10").Value = '-'
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = '-'
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = '-'
$ws.Range("N10").Value = '-'
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = '-'
$ws.Range("Q10").Value = '-'
$ws.Range("R10").Value = '1698.24 : 1'
$ws.Range("S10").Value = '-'
$ws.Range("T10").Value = '-'

# Row 11
$ws.Range("A11").Value = '2024-03-25'
$ws.Range("B11").Value = '아이엠비디엑스'
$ws.Range("C11").Value = '미래'
$ws.Range("D11").Value = '2024-03-28'
$ws.Range("E11").Value = '2024-04-03'
$ws.Range("F11").Value = 32500000
$ws.Range("G11").Value = 2500000
$ws.Range("H11").Value = '-'
$ws.Range("I11").Value = 7700
$ws.Range("J11").Value = 9900
$ws.Range("K11").Value = '-'
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = '-'
$ws.Range("N11").Value = '-'
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = '-'
$ws.Range("Q11").Value = '-'
$ws.Range("R11").Value = '2654.2 : 1'
$ws.Range("S11").Value = '-'
$ws.Range("T11").Value = '-'

# Row 12
$ws.Range("A12").Value = '2024-03-18'
$ws.Range("B12").Value = '하나32호스팩'
$ws.Range("C12").Value = '하나'
$ws.Range("D12").Value = '2024-03-21'
$ws.Range("E12").Value = '2024-03-27'
$ws.Range("F12").Value = 6000000
$ws.Range("G12").Value = 3000000
$ws.Range("H12").Value = '-'
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = '-'
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = '-'
$ws.Range("N12").Value = '-'
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = '-'
$ws.Range("Q12").Value = '-'
$ws.Range("R12").Value = '2389.8 : 1'
$ws.Range("S12").Value = '-'
$ws.Range("T12").Value = '-'

# Row 13
$ws.Range("A13").Value = '2024-03-14'
$ws.Range("B13").Value = '엔젤로보틱스'
$ws.Range("C13").Value = 'NH'
$ws.Range("D13").Value = '2024-03-19'
$ws.Range("E13").Value = '2024-03-26'
$ws.Range("F13").Value = 32000000
$ws.Range("G13").Value = 1600000
$ws.Range("H13").Value = '-'
$ws.Range("I13").Value = 11000
$ws.Range("J13").Value = 15000
$ws.Range("K13").Value = '-'
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = '-'
$ws.Range("N13").Value = '-'
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = '-'
$ws.Range("Q13").Value = '-'
$ws.Range("R13").Value = '2242.016 : 1'
$ws.Range("S13").Value = '-'
$ws.Range("T13").Value = '-'

# Row 14
$ws.Range("A14").Value = '2024-03-12'
$ws.Range("B14").Value = '삼현'
$ws.Range("C14").Value = '한국'
$ws.Range("D14").Value = '2024-03-15'
$ws.Range("E14").Value = '2024-03-21'
$ws.Range("F14").Value = 60000000
$ws.Range("G14").Value = 2000000
$ws.Range("H14").Value = '-'
$ws.Range("I14").Value = 20000
$ws.Range("J14").Value = 25000
$ws.Range("K14").Value = '-'
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = '-'
$ws.Range("N14").Value = '-'
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = '-'
$ws.Range("Q14").Value = '-'
$ws.Range("R14").Value = '1645.13 : 1'
$ws.Range("S14").Value = '-'
$ws.Range("T14").Value = '-'

# Clear the temporary Text number format back to the default "Normal"
# style now that the literal strings are safely stored as text, so the
# saved styles.xml does not pick up a stray custom number format.
$colA.Style = "Normal"
$colD.Style = "Normal"
$colE.Style = "Normal"

